# Update cryptos price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.068.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.984.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -4.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("E11").Value = "  -4.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.495.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("E15").Value = "  -5.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.132.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.987.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -5.14%  "
$ws.Range("E28").Value = "  -8.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "154.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").Value = "  -3.61%  "
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("E37").Value = "  -6.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.77%  "
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.015.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.190.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.85%  "
$ws.Range("E46").Value = "  -6.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.933"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.66%  "
$ws.Range("E49").Value = "  -4.98%  "
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("E51").Value = "  -10.97%  "
